# Update a few data values on the training schedule sheet and move the
# active selection, matching the authoring change recorded in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 data updates
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = -3
$ws.Range("H4").Value = 46

# Move the active selection from D5 to C4
$ws.Range("C4").Select() | Out-Null
